$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 3024
$ws.Range("E14").Value = 12088
$ws.Range("G14").Value = 10418
$ws.Range("H14").Value = 3861
$ws.Range("M14").Value = 1220
$ws.Range("R14").Value = 1435
$ws.Range("T14").Value = 5943
$ws.Range("U14").Value = 9669
$ws.Range("X14").Value = 5117
$ws.Range("AG14").Value = 7372
$ws.Range("AI14").Value = 93889
$ws.Range("B15").Value = 3125
$ws.Range("G15").Value = 10833
$ws.Range("H15").Value = 4362
$ws.Range("Q15").Value = 615
$ws.Range("R15").Value = 1752
$ws.Range("U15").Value = 8858
$ws.Range("AG15").Value = 7166
$ws.Range("AI15").Value = 96470
$ws.Range("B16").Value = 3456
$ws.Range("G16").Value = 11980
$ws.Range("H16").Value = 4568
$ws.Range("U16").Value = 10853
$ws.Range("X16").Value = 5592
$ws.Range("AG16").Value = 8177
$ws.Range("AI16").Value = 111413
$ws.Range("B17").Value = 3754
$ws.Range("G17").Value = 13394
$ws.Range("H17").Value = 5090
$ws.Range("K17").Value = 361
$ws.Range("R17").Value = 2243
$ws.Range("T17").Value = 7146
$ws.Range("U17").Value = 12091
$ws.Range("X17").Value = 5606
$ws.Range("Y17").Value = 3957
$ws.Range("AG17").Value = 9320
$ws.Range("AI17").Value = 121881
$ws.Range("B18").Value = 3760
$ws.Range("G18").Value = 14074
$ws.Range("H18").Value = 5538
$ws.Range("N18").Value = 2976
$ws.Range("R18").Value = 2590
$ws.Range("T18").Value = 8383
$ws.Range("U18").Value = 13985
$ws.Range("X18").Value = 6192
$ws.Range("Y18").Value = 4189
$ws.Range("AG18").Value = 10346
$ws.Range("AI18").Value = 129904
$ws.Range("G19").Value = 15366
$ws.Range("H19").Value = 6053
$ws.Range("I19").Value = 4012
$ws.Range("N19").Value = 3368
$ws.Range("T19").Value = 8976
$ws.Range("U19").Value = 15091
$ws.Range("V19").Value = 12539
$ws.Range("X19").Value = 6501
$ws.Range("Y19").Value = 4310
$ws.Range("AF19").Value = 126233
$ws.Range("AG19").Value = 11034
$ws.Range("AI19").Value = 137915
$ws.Range("B20").Value = 4538
$ws.Range("G20").Value = 16823
$ws.Range("H20").Value = 6861
$ws.Range("I20").Value = 4714
$ws.Range("Q20").Value = 853
$ws.Range("R20").Value = 2553
$ws.Range("T20").Value = 9400
$ws.Range("U20").Value = 16609
$ws.Range("V20").Value = 13776
$ws.Range("X20").Value = 6963
$ws.Range("AF20").Value = 135911
$ws.Range("AG20").Value = 11950
$ws.Range("AI20").Value = 148579
$ws.Range("B21").Value = 5197
$ws.Range("G21").Value = 18583
$ws.Range("H21").Value = 7820
$ws.Range("I21").Value = 5095
$ws.Range("P21").Value = 2324
$ws.Range("R21").Value = 2709
$ws.Range("T21").Value = 10459
$ws.Range("U21").Value = 18012
$ws.Range("V21").Value = 14793
$ws.Range("X21").Value = 8569
$ws.Range("AF21").Value = 145641
$ws.Range("AG21").Value = 13107
$ws.Range("AI21").Value = 159498
$ws.Range("B22").Value = 5876
$ws.Range("G22").Value = 18555
$ws.Range("H22").Value = 7896
$ws.Range("I22").Value = 5418
$ws.Range("P22").Value = 2220
$ws.Range("T22").Value = 11494
$ws.Range("U22").Value = 19679
$ws.Range("V22").Value = 16115
$ws.Range("X22").Value = 8966
$ws.Range("AF22").Value = 154853
$ws.Range("AG22").Value = 13719
$ws.Range("AI22").Value = 169267
$ws.Range("E23").Value = 15605
$ws.Range("G23").Value = 18751
$ws.Range("H23").Value = 8147
$ws.Range("I23").Value = 5855
$ws.Range("J23").Value = 2291
$ws.Range("R23").Value = 2904
$ws.Range("S23").Value = 5286
$ws.Range("T23").Value = 11564
$ws.Range("U23").Value = 21011
$ws.Range("V23").Value = 17092
$ws.Range("X23").Value = 8727
$ws.Range("Y23").Value = 4676
$ws.Range("AF23").Value = 164306
$ws.Range("AG23").Value = 14719
$ws.Range("AI23").Value = 179841
$ws.Range("B24").Value = 5664
$ws.Range("E24").Value = 16281
$ws.Range("G24").Value = 20203
$ws.Range("H24").Value = 8793
$ws.Range("I24").Value = 5971
$ws.Range("J24").Value = 2822
$ws.Range("L24").Value = 1124
$ws.Range("P24").Value = 2221
$ws.Range("Q24").Value = 899
$ws.Range("R24").Value = 2971
$ws.Range("T24").Value = 12401
$ws.Range("U24").Value = 22594
$ws.Range("V24").Value = 18330
$ws.Range("X24").Value = 8684
$ws.Range("AF24").Value = 174366
$ws.Range("AG24").Value = 15867
$ws.Range("AI24").Value = 191064
$ws.Range("B25").Value = 5811
$ws.Range("G25").Value = 19958
$ws.Range("H25").Value = 9068
$ws.Range("I25").Value = 6100
$ws.Range("J25").Value = 2968
$ws.Range("K25").Value = 294
$ws.Range("L25").Value = 1094
$ws.Range("R25").Value = 3179
$ws.Range("T25").Value = 13620
$ws.Range("U25").Value = 22511
$ws.Range("V25").Value = 17989
$ws.Range("X25").Value = 9295
$ws.Range("AF25").Value = 179758
$ws.Range("AI25").Value = 196486
$ws.Range("B26").Value = 6829
$ws.Range("D26").Value = 24837
$ws.Range("E26").Value = 22305
$ws.Range("G26").Value = 19827
$ws.Range("H26").Value = 8987
$ws.Range("I26").Value = 6253
$ws.Range("J26").Value = 2734
$ws.Range("L26").Value = 1097
$ws.Range("P26").Value = 2377
$ws.Range("Q26").Value = 756
$ws.Range("T26").Value = 11966
$ws.Range("U26").Value = 21728
$ws.Range("V26").Value = 18543
$ws.Range("X26").Value = 8295
$ws.Range("Y26").Value = 4171
$ws.Range("AF26").Value = 183526
$ws.Range("AG26").Value = 15562
$ws.Range("AI26").Value = 199990
